{"js": "// Replace each old three-digit-by-one-digit multiplication equation with\n// the new equation, matching the document order the diff specifies.\n// Each source string is unique in the document, so an exact, case-\n// sensitive search-and-replace on the whole body is sufficient.\n\nconst replacements = [\n  [\"746\u00d74=2984\", \"290\u00d72=580\"],\n  [\"653\u00d72=1306\", \"111\u00d76=666\"],\n  [\"177\u00d77=1239\", \"839\u00d77=5873\"],\n  [\"624\u00d79=5616\", \"722\u00d79=6498\"],\n  [\"806\u00d73=2418\", \"465\u00d78=3720\"],\n  [\"200\u00d77=1400\", \"743\u00d76=4458\"],\n  [\"609\u00d73=1827\", \"968\u00d73=2904\"],\n  [\"885\u00d75=4425\", \"794\u00d79=7146\"],\n  [\"927\u00d72=1854\", \"628\u00d73=1884\"],\n  [\"355\u00d76=2130\", \"146\u00d76=876\"],\n  [\"559\u00d78=4472\", \"481\u00d74=1924\"],\n  [\"937\u00d78=7496\", \"211\u00d77=1477\"],\n  [\"209\u00d79=1881\", \"414\u00d76=2484\"],\n  [\"618\u00d76=3708\", \"278\u00d79=2502\"],\n  [\"795\u00d72=1590\", \"785\u00d73=2355\"],\n  [\"578\u00d72=1156\", \"766\u00d79=6894\"],\n  [\"924\u00d76=5544\", \"459\u00d78=3672\"],\n  [\"180\u00d73=540\", \"120\u00d75=600\"],\n  [\"878\u00d77=6146\", \"182\u00d73=546\"],\n  [\"262\u00d74=1048\", \"118\u00d72=236\"],\n  [\"988\u00d72=1976\", \"309\u00d74=1236\"],\n  [\"483\u00d74=1932\", \"228\u00d75=1140\"],\n  [\"224\u00d73=672\", \"376\u00d76=2256\"],\n  [\"365\u00d73=1095\", \"281\u00d77=1967\"],\n  [\"331\u00d78=2648\", \"329\u00d77=2303\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-by-one-digit multiplication equation with\n# the new equation, using Word's Find/Replace on the whole document\n# content. Every source string is unique in the document, so a single\n# wdReplaceAll pass per pair is sufficient and safely idempotent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"746\u00d74=2984\"; New = \"290\u00d72=580\" },\n    @{ Old = \"653\u00d72=1306\"; New = \"111\u00d76=666\" },\n    @{ Old = \"177\u00d77=1239\"; New = \"839\u00d77=5873\" },\n    @{ Old = \"624\u00d79=5616\"; New = \"722\u00d79=6498\" },\n    @{ Old = \"806\u00d73=2418\"; New = \"465\u00d78=3720\" },\n    @{ Old = \"200\u00d77=1400\"; New = \"743\u00d76=4458\" },\n    @{ Old = \"609\u00d73=1827\"; New = \"968\u00d73=2904\" },\n    @{ Old = \"885\u00d75=4425\"; New = \"794\u00d79=7146\" },\n    @{ Old = \"927\u00d72=1854\"; New = \"628\u00d73=1884\" },\n    @{ Old = \"355\u00d76=2130\"; New = \"146\u00d76=876\" },\n    @{ Old = \"559\u00d78=4472\"; New = \"481\u00d74=1924\" },\n    @{ Old = \"937\u00d78=7496\"; New = \"211\u00d77=1477\" },\n    @{ Old = \"209\u00d79=1881\"; New = \"414\u00d76=2484\" },\n    @{ Old = \"618\u00d76=3708\"; New = \"278\u00d79=2502\" },\n    @{ Old = \"795\u00d72=1590\"; New = \"785\u00d73=2355\" },\n    @{ Old = \"578\u00d72=1156\"; New = \"766\u00d79=6894\" },\n    @{ Old = \"924\u00d76=5544\"; New = \"459\u00d78=3672\" },\n    @{ Old = \"180\u00d73=540\"; New = \"120\u00d75=600\" },\n    @{ Old = \"878\u00d77=6146\"; New = \"182\u00d73=546\" },\n    @{ Old = \"262\u00d74=1048\"; New = \"118\u00d72=236\" },\n    @{ Old = \"988\u00d72=1976\"; New = \"309\u00d74=1236\" },\n    @{ Old = \"483\u00d74=1932\"; New = \"228\u00d75=1140\" },\n    @{ Old = \"224\u00d73=672\"; New = \"376\u00d76=2256\" },\n    @{ Old = \"365\u00d73=1095\"; New = \"281\u00d77=1967\" },\n    @{ Old = \"331\u00d78=2648\"; New = \"329\u00d77=2303\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
